$wb = $excel.ActiveWorkbook

# --- PIR sheet: append rows 594-606 ---
$ws = $wb.Worksheets.Item("PIR")
$rng = $ws.Range("A594:F606")
$rng.NumberFormat = "@"
$ws.Cells.Item(594, 1).Value = "2026-02-06"
$ws.Cells.Item(594, 2).Value = "10:29:18"
$ws.Cells.Item(594, 3).Value = "10:00"
$ws.Cells.Item(594, 4).Value = "Bathroom"
$ws.Cells.Item(594, 5).Value = "No Motion"
$ws.Cells.Item(594, 6).Value = "Inactive"
$ws.Cells.Item(595, 1).Value = "2026-02-06"
$ws.Cells.Item(595, 2).Value = "10:29:22"
$ws.Cells.Item(595, 3).Value = "10:00"
$ws.Cells.Item(595, 4).Value = "Bathroom"
$ws.Cells.Item(595, 5).Value = "No Motion"
$ws.Cells.Item(595, 6).Value = "Inactive"
$ws.Cells.Item(596, 1).Value = "2026-02-06"
$ws.Cells.Item(596, 2).Value = "10:29:26"
$ws.Cells.Item(596, 3).Value = "10:00"
$ws.Cells.Item(596, 4).Value = "Bathroom"
$ws.Cells.Item(596, 5).Value = "No Motion"
$ws.Cells.Item(596, 6).Value = "Inactive"
$ws.Cells.Item(597, 1).Value = "2026-02-06"
$ws.Cells.Item(597, 2).Value = "10:29:30"
$ws.Cells.Item(597, 3).Value = "10:00"
$ws.Cells.Item(597, 4).Value = "Bathroom"
$ws.Cells.Item(597, 5).Value = "No Motion"
$ws.Cells.Item(597, 6).Value = "Inactive"
$ws.Cells.Item(598, 1).Value = "2026-02-06"
$ws.Cells.Item(598, 2).Value = "10:29:33"
$ws.Cells.Item(598, 3).Value = "10:00"
$ws.Cells.Item(598, 4).Value = "Bathroom"
$ws.Cells.Item(598, 5).Value = "No Motion"
$ws.Cells.Item(598, 6).Value = "Inactive"
$ws.Cells.Item(599, 1).Value = "2026-02-06"
$ws.Cells.Item(599, 2).Value = "10:29:38"
$ws.Cells.Item(599, 3).Value = "10:00"
$ws.Cells.Item(599, 4).Value = "Bathroom"
$ws.Cells.Item(599, 5).Value = "No Motion"
$ws.Cells.Item(599, 6).Value = "Inactive"
$ws.Cells.Item(600, 1).Value = "2026-02-06"
$ws.Cells.Item(600, 2).Value = "10:29:43"
$ws.Cells.Item(600, 3).Value = "10:00"
$ws.Cells.Item(600, 4).Value = "Bathroom"
$ws.Cells.Item(600, 5).Value = "No Motion"
$ws.Cells.Item(600, 6).Value = "Inactive"
$ws.Cells.Item(601, 1).Value = "2026-02-06"
$ws.Cells.Item(601, 2).Value = "10:29:48"
$ws.Cells.Item(601, 3).Value = "10:00"
$ws.Cells.Item(601, 4).Value = "Bathroom"
$ws.Cells.Item(601, 5).Value = "No Motion"
$ws.Cells.Item(601, 6).Value = "Inactive"
$ws.Cells.Item(602, 1).Value = "2026-02-06"
$ws.Cells.Item(602, 2).Value = "10:29:53"
$ws.Cells.Item(602, 3).Value = "10:00"
$ws.Cells.Item(602, 4).Value = "Bathroom"
$ws.Cells.Item(602, 5).Value = "No Motion"
$ws.Cells.Item(602, 6).Value = "Inactive"
$ws.Cells.Item(603, 1).Value = "2026-02-06"
$ws.Cells.Item(603, 2).Value = "10:29:58"
$ws.Cells.Item(603, 3).Value = "10:00"
$ws.Cells.Item(603, 4).Value = "Bathroom"
$ws.Cells.Item(603, 5).Value = "No Motion"
$ws.Cells.Item(603, 6).Value = "Inactive"
$ws.Cells.Item(604, 1).Value = "2026-02-06"
$ws.Cells.Item(604, 2).Value = "10:30:03"
$ws.Cells.Item(604, 3).Value = "10:00"
$ws.Cells.Item(604, 4).Value = "Bathroom"
$ws.Cells.Item(604, 5).Value = "No Motion"
$ws.Cells.Item(604, 6).Value = "Inactive"
$ws.Cells.Item(605, 1).Value = "2026-02-06"
$ws.Cells.Item(605, 2).Value = "10:30:08"
$ws.Cells.Item(605, 3).Value = "10:00"
$ws.Cells.Item(605, 4).Value = "Bathroom"
$ws.Cells.Item(605, 5).Value = "No Motion"
$ws.Cells.Item(605, 6).Value = "Inactive"
$ws.Cells.Item(606, 1).Value = "2026-02-06"
$ws.Cells.Item(606, 2).Value = "10:30:13"
$ws.Cells.Item(606, 3).Value = "10:00"
$ws.Cells.Item(606, 4).Value = "Bathroom"
$ws.Cells.Item(606, 5).Value = "No Motion"
$ws.Cells.Item(606, 6).Value = "Inactive"
$rng.ClearFormats()

# --- Humidity sheet: append rows 424-435 ---
$ws = $wb.Worksheets.Item("Humidity")
$rng = $ws.Range("A424:F435")
$rng.NumberFormat = "@"
$ws.Cells.Item(424, 1).Value = "2026-02-06"
$ws.Cells.Item(424, 2).Value = "10:29:16"
$ws.Cells.Item(424, 3).Value = "10:00"
$ws.Cells.Item(424, 4).Value = "Bathroom"
$ws.Cells.Item(424, 5).Value = "67.7%"
$ws.Cells.Item(424, 6).Value = "Active"
$ws.Cells.Item(425, 1).Value = "2026-02-06"
$ws.Cells.Item(425, 2).Value = "10:29:20"
$ws.Cells.Item(425, 3).Value = "10:00"
$ws.Cells.Item(425, 4).Value = "Bathroom"
$ws.Cells.Item(425, 5).Value = "67.7%"
$ws.Cells.Item(425, 6).Value = "Active"
$ws.Cells.Item(426, 1).Value = "2026-02-06"
$ws.Cells.Item(426, 2).Value = "10:29:24"
$ws.Cells.Item(426, 3).Value = "10:00"
$ws.Cells.Item(426, 4).Value = "Bathroom"
$ws.Cells.Item(426, 5).Value = "66.7%"
$ws.Cells.Item(426, 6).Value = "Active"
$ws.Cells.Item(427, 1).Value = "2026-02-06"
$ws.Cells.Item(427, 2).Value = "10:29:28"
$ws.Cells.Item(427, 3).Value = "10:00"
$ws.Cells.Item(427, 4).Value = "Bathroom"
$ws.Cells.Item(427, 5).Value = "67.6%"
$ws.Cells.Item(427, 6).Value = "Active"
$ws.Cells.Item(428, 1).Value = "2026-02-06"
$ws.Cells.Item(428, 2).Value = "10:29:34"
$ws.Cells.Item(428, 3).Value = "10:00"
$ws.Cells.Item(428, 4).Value = "Bathroom"
$ws.Cells.Item(428, 5).Value = "66.2%"
$ws.Cells.Item(428, 6).Value = "Active"
$ws.Cells.Item(429, 1).Value = "2026-02-06"
$ws.Cells.Item(429, 2).Value = "10:29:39"
$ws.Cells.Item(429, 3).Value = "10:00"
$ws.Cells.Item(429, 4).Value = "Bathroom"
$ws.Cells.Item(429, 5).Value = "66.5%"
$ws.Cells.Item(429, 6).Value = "Active"
$ws.Cells.Item(430, 1).Value = "2026-02-06"
$ws.Cells.Item(430, 2).Value = "10:29:44"
$ws.Cells.Item(430, 3).Value = "10:00"
$ws.Cells.Item(430, 4).Value = "Bathroom"
$ws.Cells.Item(430, 5).Value = "67.5%"
$ws.Cells.Item(430, 6).Value = "Active"
$ws.Cells.Item(431, 1).Value = "2026-02-06"
$ws.Cells.Item(431, 2).Value = "10:29:49"
$ws.Cells.Item(431, 3).Value = "10:00"
$ws.Cells.Item(431, 4).Value = "Bathroom"
$ws.Cells.Item(431, 5).Value = "67.5%"
$ws.Cells.Item(431, 6).Value = "Active"
$ws.Cells.Item(432, 1).Value = "2026-02-06"
$ws.Cells.Item(432, 2).Value = "10:29:55"
$ws.Cells.Item(432, 3).Value = "10:00"
$ws.Cells.Item(432, 4).Value = "Bathroom"
$ws.Cells.Item(432, 5).Value = "67.5%"
$ws.Cells.Item(432, 6).Value = "Active"
$ws.Cells.Item(433, 1).Value = "2026-02-06"
$ws.Cells.Item(433, 2).Value = "10:30:00"
$ws.Cells.Item(433, 3).Value = "10:00"
$ws.Cells.Item(433, 4).Value = "Bathroom"
$ws.Cells.Item(433, 5).Value = "67.5%"
$ws.Cells.Item(433, 6).Value = "Active"
$ws.Cells.Item(434, 1).Value = "2026-02-06"
$ws.Cells.Item(434, 2).Value = "10:30:10"
$ws.Cells.Item(434, 3).Value = "10:00"
$ws.Cells.Item(434, 4).Value = "Bathroom"
$ws.Cells.Item(434, 5).Value = "66.5%"
$ws.Cells.Item(434, 6).Value = "Active"
$ws.Cells.Item(435, 1).Value = "2026-02-06"
$ws.Cells.Item(435, 2).Value = "10:30:15"
$ws.Cells.Item(435, 3).Value = "10:00"
$ws.Cells.Item(435, 4).Value = "Bathroom"
$ws.Cells.Item(435, 5).Value = "67.6%"
$ws.Cells.Item(435, 6).Value = "Active"
$rng.ClearFormats()

# --- Temperature sheet: append rows 424-434 ---
$ws = $wb.Worksheets.Item("Temperature")
$rng = $ws.Range("A424:F434")
$rng.NumberFormat = "@"
$ws.Cells.Item(424, 1).Value = "2026-02-06"
$ws.Cells.Item(424, 2).Value = "10:29:17"
$ws.Cells.Item(424, 3).Value = "10:00"
$ws.Cells.Item(424, 4).Value = "Bathroom"
$ws.Cells.Item(424, 5).Value = "28.5C"
$ws.Cells.Item(424, 6).Value = "Active"
$ws.Cells.Item(425, 1).Value = "2026-02-06"
$ws.Cells.Item(425, 2).Value = "10:29:21"
$ws.Cells.Item(425, 3).Value = "10:00"
$ws.Cells.Item(425, 4).Value = "Bathroom"
$ws.Cells.Item(425, 5).Value = "28.5C"
$ws.Cells.Item(425, 6).Value = "Active"
$ws.Cells.Item(426, 1).Value = "2026-02-06"
$ws.Cells.Item(426, 2).Value = "10:29:25"
$ws.Cells.Item(426, 3).Value = "10:00"
$ws.Cells.Item(426, 4).Value = "Bathroom"
$ws.Cells.Item(426, 5).Value = "28.4C"
$ws.Cells.Item(426, 6).Value = "Active"
$ws.Cells.Item(427, 1).Value = "2026-02-06"
$ws.Cells.Item(427, 2).Value = "10:29:29"
$ws.Cells.Item(427, 3).Value = "10:00"
$ws.Cells.Item(427, 4).Value = "Bathroom"
$ws.Cells.Item(427, 5).Value = "28.4C"
$ws.Cells.Item(427, 6).Value = "Active"
$ws.Cells.Item(428, 1).Value = "2026-02-06"
$ws.Cells.Item(428, 2).Value = "10:29:36"
$ws.Cells.Item(428, 3).Value = "10:00"
$ws.Cells.Item(428, 4).Value = "Bathroom"
$ws.Cells.Item(428, 5).Value = "28.5C"
$ws.Cells.Item(428, 6).Value = "Active"
$ws.Cells.Item(429, 1).Value = "2026-02-06"
$ws.Cells.Item(429, 2).Value = "10:29:41"
$ws.Cells.Item(429, 3).Value = "10:00"
$ws.Cells.Item(429, 4).Value = "Bathroom"
$ws.Cells.Item(429, 5).Value = "28.4C"
$ws.Cells.Item(429, 6).Value = "Active"
$ws.Cells.Item(430, 1).Value = "2026-02-06"
$ws.Cells.Item(430, 2).Value = "10:29:46"
$ws.Cells.Item(430, 3).Value = "10:00"
$ws.Cells.Item(430, 4).Value = "Bathroom"
$ws.Cells.Item(430, 5).Value = "28.4C"
$ws.Cells.Item(430, 6).Value = "Active"
$ws.Cells.Item(431, 1).Value = "2026-02-06"
$ws.Cells.Item(431, 2).Value = "10:29:51"
$ws.Cells.Item(431, 3).Value = "10:00"
$ws.Cells.Item(431, 4).Value = "Bathroom"
$ws.Cells.Item(431, 5).Value = "28.4C"
$ws.Cells.Item(431, 6).Value = "Active"
$ws.Cells.Item(432, 1).Value = "2026-02-06"
$ws.Cells.Item(432, 2).Value = "10:29:56"
$ws.Cells.Item(432, 3).Value = "10:00"
$ws.Cells.Item(432, 4).Value = "Bathroom"
$ws.Cells.Item(432, 5).Value = "28.4C"
$ws.Cells.Item(432, 6).Value = "Active"
$ws.Cells.Item(433, 1).Value = "2026-02-06"
$ws.Cells.Item(433, 2).Value = "10:30:01"
$ws.Cells.Item(433, 3).Value = "10:00"
$ws.Cells.Item(433, 4).Value = "Bathroom"
$ws.Cells.Item(433, 5).Value = "28.4C"
$ws.Cells.Item(433, 6).Value = "Active"
$ws.Cells.Item(434, 1).Value = "2026-02-06"
$ws.Cells.Item(434, 2).Value = "10:30:11"
$ws.Cells.Item(434, 3).Value = "10:00"
$ws.Cells.Item(434, 4).Value = "Bathroom"
$ws.Cells.Item(434, 5).Value = "28.3C"
$ws.Cells.Item(434, 6).Value = "Active"
$rng.ClearFormats()
